# Applies updated "想去人数" (F column) figures across sheets, matching
# the published gh-pages data snapshot at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 285
$ws.Range("F4").Value = 1828
$ws.Range("F10").Value = 3612
$ws.Range("F11").Value = 149
$ws.Range("F12").Value = 97
$ws.Range("F14").Value = 57
$ws.Range("F15").Value = 65
$ws.Range("F16").Value = 627
$ws.Range("F17").Value = 124
$ws.Range("F18").Value = 805
$ws.Range("F19").Value = 27
$ws.Range("F21").Value = 137
$ws.Range("F25").Value = 2904
$ws.Range("F26").Value = 5340
$ws.Range("F29").Value = 488
$ws.Range("F30").Value = 20
$ws.Range("F31").Value = 3126
$ws.Range("F32").Value = 312
$ws.Range("F33").Value = 2316
$ws.Range("F34").Value = 570
$ws.Range("F36").Value = 95
$ws.Range("F37").Value = 148
$ws.Range("F38").Value = 201
$ws.Range("F39").Value = 318
$ws.Range("F40").Value = 62
$ws.Range("F41").Value = 477
$ws.Range("F42").Value = 829
$ws.Range("F43").Value = 37

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 78

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 285
$ws.Range("F4").Value = 1828
$ws.Range("F10").Value = 3612
$ws.Range("F11").Value = 149
$ws.Range("F12").Value = 97
$ws.Range("F14").Value = 78
$ws.Range("F15").Value = 57
$ws.Range("F16").Value = 65
$ws.Range("F17").Value = 627
$ws.Range("F18").Value = 124
$ws.Range("F19").Value = 805
$ws.Range("F20").Value = 27
$ws.Range("F22").Value = 137
$ws.Range("F26").Value = 2904
$ws.Range("F27").Value = 5340
$ws.Range("F30").Value = 488
$ws.Range("F31").Value = 20
$ws.Range("F32").Value = 3126
$ws.Range("F33").Value = 312
$ws.Range("F34").Value = 2316
$ws.Range("F35").Value = 570
$ws.Range("F37").Value = 95
$ws.Range("F38").Value = 148
$ws.Range("F39").Value = 201
$ws.Range("F40").Value = 318
$ws.Range("F41").Value = 62
$ws.Range("F42").Value = 477
$ws.Range("F43").Value = 829
$ws.Range("F44").Value = 37
